$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF4").Value = 0.681
$ws.Range("AF5").Value = 1
$ws.Range("AF6").Value = 0.8100000000000001
$ws.Range("AF7").Value = 0.914
$ws.Range("AF8").Value = 0.9379999999999999
$ws.Range("AF9").Value = 0.833
$ws.Range("AF10").Value = 1
$ws.Range("AF11").Value = 1
$ws.Range("AF12").Value = 1.167
$ws.Range("AF13").Value = 2
